$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Models sheet: add STRATIFICATION column (E) and update its selection
# ---------------------------------------------------------------------------
$wsModels = $wb.Worksheets.Item("Models")
$wsModels.Range("E1").Value = "STRATIFICATION"
$wsModels.Columns.Item(5).ColumnWidth = 35.5

# ---------------------------------------------------------------------------
# SubjectMetabolites sheet: just move the selection
# ---------------------------------------------------------------------------
$wsSubjMeta = $wb.Worksheets.Item("SubjectMetabolites")
$wsSubjMeta.Range("C24").Select() | Out-Null

# ---------------------------------------------------------------------------
# Models sheet selection (do this after SubjectMetabolites so tab order works
# itself out once SubjectData is activated last)
# ---------------------------------------------------------------------------
$wsModels.Range("D10").Select() | Out-Null

# ---------------------------------------------------------------------------
# SubjectData sheet: insert a new SITE column (B) with stratification codes
# ---------------------------------------------------------------------------
$wsSubjData = $wb.Worksheets.Item("SubjectData")
$wsSubjData.Columns("B:B").Insert()
$wsSubjData.Columns.Item(2).ColumnWidth = 13.45

$wsSubjData.Range("B1").Value = "SITE"

$siteValues = @(1,1,1,2,3,1,1,2,2,2,3,3,3,3,2,2,3,2,1)
for ($i = 0; $i -lt $siteValues.Length; $i++) {
    $row = $i + 2
    $wsSubjData.Cells.Item($row, 2).Value = $siteValues[$i]
}

# SubjectData becomes the active sheet/tab, with the new selection
$wsSubjData.Range("B10").Select() | Out-Null
